$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new value next to "valor hora técnica:$50/h" (row 21), matching the
# style already used by the other labor-cost rows below it (B23/B24/B26).
$ws.Range("B21").Value = 1800
$ws.Range("B21").HorizontalAlignment = $ws.Range("B23").HorizontalAlignment
$ws.Range("B21").NumberFormat = $ws.Range("B23").NumberFormat

# Update "valor backup:" amount.
$ws.Range("B24").Value = 1200

# Update "valor formatação:" amount and bring its format in line with the
# sibling rows (was General/left, now #,##0/left like B23/B24).
$ws.Range("B25").NumberFormat = $ws.Range("B23").NumberFormat
$ws.Range("B25").Value = 1000

# Scroll the view down and move the selection, matching where the user left
# off after entering the new figures.
$ws.Range("B28").Select()
$ws.Application.ActiveWindow.ScrollRow = 10
